$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 114, pushing existing rows 114:206 down to 115:207.
$ws.Rows.Item(114).Insert()

# Populate the newly inserted row 114 with the new weekly data point
# (same product/quality/origin as the row that used to occupy 114, but a new
# sampling date and updated min/max/weighted-average prices).
$ws.Range("A114").Value = 9
$ws.Range("B114").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C114").Value = "Metropolitana"
$ws.Range("D114").Value = 44483
$ws.Range("E114").Value = 13
$ws.Range("F114").Value = 100112028
$ws.Range("G114").Value = "Sandia"
$ws.Range("H114").Value = "Sin especificar"
$ws.Range("I114").Value = "Primera"
$ws.Range("J114").Value = 250
$ws.Range("K114").Value = 800
$ws.Range("L114").Value = 900
$ws.Range("M114").Value = 850
$ws.Range("N114").Value = "$/kilo (volumen en unidades)"
$ws.Range("O114").Value = "Perú"
$ws.Range("P114").Value = 850
$ws.Range("Q114").Value = 1
$ws.Range("R114").Value = "Hortaliza"
